$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.794.68"
$ws.Range("E2").Value = "  +4.02%  "
$ws.Range("D3").Value = "2.424.72"
$ws.Range("E3").Value = "  +2.42%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.34"
$ws.Range("E5").Value = "  +3.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.91"
$ws.Range("E6").Value = "  +7.00%  "
$ws.Range("E7").Value = "  +2.65%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.533"
$ws.Range("E9").Value = "  +10.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.49"
$ws.Range("E10").Value = "  +3.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0801"
$ws.Range("E11").Value = "  +1.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "18.82"
$ws.Range("E12").Value = "  +1.26%  "
$ws.Range("E13").Value = "  -1.32%  "
$ws.Range("E14").Value = "  +3.21%  "
$ws.Range("D15").Value = "2.802.69"
$ws.Range("D16").Value = "2.412.56"
$ws.Range("E16").Value = "  +2.40%  "
$ws.Range("E17").Value = "  +4.46%  "
$ws.Range("D18").Value = "44.629.23"
$ws.Range("E18").Value = "  +3.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.31"
$ws.Range("E19").Value = "  +2.75%  "
$ws.Range("E20").Value = "  +1.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.79"
$ws.Range("E22").Value = "  +1.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "242.33"
$ws.Range("E23").Value = "  +2.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.27"
$ws.Range("E24").Value = "  +4.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.49"
$ws.Range("E25").Value = "  +2.34%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.26"
$ws.Range("E27").Value = "  +2.98%  "
$ws.Range("E28").Value = "  -3.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.48"
$ws.Range("E29").Value = "  +1.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.75"
$ws.Range("E30").Value = "  +4.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "48.60"
$ws.Range("E31").Value = "  +1.60%  "
$ws.Range("E32").Value = "  +18.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.52"
$ws.Range("E33").Value = "  +11.34%  "
$ws.Range("E34").Value = "  +3.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0775"
$ws.Range("E35").Value = "  +6.27%  "
$ws.Range("E37").Value = "  +3.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.47"
$ws.Range("E38").Value = "  +3.23%  "
$ws.Range("E39").Value = "  +0.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "122.77"
$ws.Range("E40").Value = "  -4.00%  "
$ws.Range("E41").Value = "  +1.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.20"
$ws.Range("E42").Value = "  -3.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.04"
$ws.Range("E43").Value = "  +1.14%  "
$ws.Range("E44").Value = "  +4.35%  "
$ws.Range("D45").Value = "1.941.56"
$ws.Range("E45").Value = "  +0.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.12"
$ws.Range("E46").Value = "  -1.12%  "
$ws.Range("E47").Value = "  +8.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.39"
$ws.Range("E48").Value = "  +0.95%  "
$ws.Range("E49").Value = "  +14.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.38"
$ws.Range("E50").Value = "  +4.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.06"
